$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timestamps")

# --- New data rows (8, 9, 10) -------------------------------------------
# Row 8: 2022-10-10, D. Hoyer, Hours column F = "?"
$ws.Range("A8").Value = (Get-Date -Year 2022 -Month 10 -Day 10).Date
$ws.Range("B8").Value = "D. Hoyer"
$ws.Range("F8").Value = "?"

# Row 9: 2022-10-10, L. Roth, Hours column F = 1
$ws.Range("A9").Value = (Get-Date -Year 2022 -Month 10 -Day 10).Date
$ws.Range("B9").Value = "L. Roth"
$ws.Range("F9").Value = 1

# Row 10: 2022-10-11, L. Roth, Task "Structs and header files", Hours column F = 0.5
$ws.Range("A10").Value = (Get-Date -Year 2022 -Month 10 -Day 11).Date
$ws.Range("B10").Value = "L. Roth"
$ws.Range("C10").Value = "Structs and header files"
$ws.Range("F10").Value = 0.5

# --- Selection moved to C13:E13 -----------------------------------------
$null = $ws.Range("C13:E13").Select()
